$wb = $excel.ActiveWorkbook

# Sheet 1
$ws = $wb.Worksheets.Item(1)
$ws.Range('H2').Value = 155
$ws.Range('I2').Value = 155
$ws.Range('K2').Value = 155
$ws.Range('M2').Value = -42
$ws.Range('H32').Value = 1810.2941
$ws.Range('J32').Value = 2046
$ws.Range('L32').Value = 2046
$ws.Range('N32').Value = -2698
$ws.Range('H38').Value = 1049.5
$ws.Range('I38').Value = 259.4
$ws.Range('K38').Value = 778.1999999999999
$ws.Range('M38').Value = -406.1999999999999
$ws.Range('H39').Value = 1079.7
$ws.Range('I39').Value = 154
$ws.Range('J39').Value = 1476.4286
$ws.Range('K39').Value = 462
$ws.Range('L39').Value = 4429.2858
$ws.Range('M39').Value = -166
$ws.Range('N39').Value = -5021.2858
$ws.Range('H43').Value = 1702.4849
$ws.Range('J43').Value = 1463.1428
$ws.Range('L43').Value = 1463.1428
$ws.Range('N43').Value = -1601.1428
$ws.Range('H44').Value = 1111111
$ws.Range('J44').Value = 0
$ws.Range('L44').Value = 0
$ws.Range('N44').ClearContents()
$ws.Range('H57').Value = 60995
$ws.Range('I57').Value = 61000
$ws.Range('J57').Value = 60990
$ws.Range('K57').Value = 183000
$ws.Range('L57').Value = 182970
$ws.Range('M57').Value = -182501
$ws.Range('N57').Value = -183968
$ws.Range('H74').Value = 13233.6
$ws.Range('I74').Value = 11926.223
$ws.Range('K74').Value = 11926.223
$ws.Range('M74').Value = -10990.223
$ws.Range('H77').Value = 13233.6
$ws.Range('I77').Value = 11926.223
$ws.Range('K77').Value = 59631.115
$ws.Range('M77').Value = -54951.115
$ws.Range('H134').Value = 88925.42999999999
$ws.Range('J134').Value = 88925.42999999999
$ws.Range('L134').Value = 88925.42999999999
$ws.Range('N134').Value = -99065.42999999999
$ws.Range('H135').Value = 2399.2307
$ws.Range('I135').Value = 2169
$ws.Range('K135').Value = 19521
$ws.Range('M135').Value = -16986
$ws.Range('H136').Value = 49958.332
$ws.Range('J136').Value = 49958.332
$ws.Range('L136').Value = 49958.332
$ws.Range('N136').Value = -60158.332
$ws.Range('H138').Value = 5884.3784
$ws.Range('I138').Value = 3608.3333
$ws.Range('J138').Value = 6786.2075
$ws.Range('K138').Value = 10824.9999
$ws.Range('L138').Value = 20358.6225
$ws.Range('M138').Value = -5684.999899999999
$ws.Range('N138').Value = -30638.6225
$ws.Range('H140').Value = 57999.8
$ws.Range('J140').Value = 57999.8
$ws.Range('L140').Value = 57999.8
$ws.Range('N140').Value = -68359.8
$ws.Range('H141').Value = 8320.4
$ws.Range('I141').Value = 8833
$ws.Range('J141').Value = 7551.5
$ws.Range('K141').Value = 26499
$ws.Range('L141').Value = 22654.5
$ws.Range('M141').Value = -21319
$ws.Range('N141').Value = -33014.5

# Sheet 2
$ws = $wb.Worksheets.Item(2)
$ws.Range('H45').Value = 2095.2104
$ws.Range('I45').Value = 1561.9231
$ws.Range('K45').Value = 1561.9231
$ws.Range('M45').Value = -1184.9231
$ws.Range('H47').Value = 30000
$ws.Range('J47').Value = 30000
$ws.Range('L47').Value = 30000
$ws.Range('N47').Value = -31450
$ws.Range('H61').Value = 3742.0908
$ws.Range('I61').Value = 3654.04
$ws.Range('K61').Value = 3654.04
$ws.Range('M61').Value = -3442.04
$ws.Range('H110').Value = 174337.11
$ws.Range('I110').Value = 219400.78
$ws.Range('K110').Value = 219400.78
$ws.Range('M110').Value = -217355.78
$ws.Range('H132').Value = 4211.6113
$ws.Range('I132').Value = 3544.7576
$ws.Range('K132').Value = 10634.2728
$ws.Range('M132').Value = -8104.272799999999
$ws.Range('H136').Value = 3742.0908
$ws.Range('I136').Value = 3654.04
$ws.Range('K136').Value = 10962.12
$ws.Range('M136').Value = -8412.119999999999

# Sheet 3
$ws = $wb.Worksheets.Item(3)
$ws.Range('H60').Value = 88926.336
$ws.Range('J60').Value = 88926.336
$ws.Range('L60').Value = 88926.336
$ws.Range('N60').Value = -90124.336
$ws.Range('H68').Value = 0
$ws.Range('J68').Value = 0
$ws.Range('L68').Value = 0
$ws.Range('N68').ClearContents()
$ws.Range('H71').Value = 0
$ws.Range('J71').Value = 0
$ws.Range('L71').Value = 0
$ws.Range('N71').ClearContents()
$ws.Range('H107').Value = 358584.03
$ws.Range('I107').Value = 1237.05
$ws.Range('J107').Value = 1251951.5
$ws.Range('K107').Value = 1237.05
$ws.Range('L107').Value = 1251951.5
$ws.Range('M107').Value = 682.95
$ws.Range('N107').Value = -1255791.5

# Sheet 4
$ws = $wb.Worksheets.Item(4)
$ws.Range('H22').Value = 586.46155
$ws.Range('I22').Value = 496.2857
$ws.Range('J22').Value = 691.6667
$ws.Range('K22').Value = 496.2857
$ws.Range('L22').Value = 691.6667
$ws.Range('M22').Value = -146.2857
$ws.Range('N22').Value = -1391.6667
$ws.Range('H52').Value = 78203.89
$ws.Range('J52').Value = 78677.86
$ws.Range('L52').Value = 78677.86
$ws.Range('N52').Value = -79265.86
$ws.Range('H58').Value = 2513.7727
$ws.Range('I58').Value = 2876.2666
$ws.Range('J58').Value = 1737
$ws.Range('K58').Value = 2876.2666
$ws.Range('L58').Value = 1737
$ws.Range('M58').Value = -2673.2666
$ws.Range('N58').Value = -2143
$ws.Range('H132').Value = 1750.225
$ws.Range('I132').Value = 1513.4286
$ws.Range('K132').Value = 4540.2858
$ws.Range('M132').Value = -2010.2858
$ws.Range('H134').Value = 337090.34
$ws.Range('I134').Value = 3779.4814
$ws.Range('K134').Value = 11338.4442
$ws.Range('M134').Value = -8803.4442
$ws.Range('H135').Value = 50823.08
$ws.Range('I135').Value = 60709
$ws.Range('J135').Value = 49999.25
$ws.Range('K135').Value = 60709
$ws.Range('L135').Value = 49999.25
$ws.Range('M135').Value = -55639
$ws.Range('N135').Value = -60139.25
$ws.Range('H136').Value = 2513.7727
$ws.Range('I136').Value = 2876.2666
$ws.Range('J136').Value = 1737
$ws.Range('K136').Value = 8628.799800000001
$ws.Range('L136').Value = 5211
$ws.Range('M136').Value = -6078.799800000001
$ws.Range('N136').Value = -10311
$ws.Range('H137').Value = 50000
$ws.Range('J137').Value = 50000
$ws.Range('L137').Value = 50000
$ws.Range('N137').Value = -60200
$ws.Range('H139').Value = 245994.67
$ws.Range('J139').Value = 245994.67
$ws.Range('L139').Value = 245994.67
$ws.Range('N139').Value = -256274.67

# Sheet 5
$ws = $wb.Worksheets.Item(5)
$ws.Range('H7').Value = 0
$ws.Range('I7').Value = 0
$ws.Range('K7').Value = 0
$ws.Range('M7').ClearContents()
$ws.Range('H39').Value = 10885.429
$ws.Range('J39').Value = 13249.9375
$ws.Range('L39').Value = 39749.8125
$ws.Range('N39').Value = -40337.8125

# Sheet 6
$ws = $wb.Worksheets.Item(6)
$ws.Range('H122').Value = 3109.125
$ws.Range('I122').Value = 2227.4736
$ws.Range('K122').Value = 6682.4208
$ws.Range('M122').Value = -4232.4208
$ws.Range('H132').Value = 33601.06
$ws.Range('I132').Value = 6079.577
$ws.Range('J132').Value = 113107.555
$ws.Range('K132').Value = 18238.731
$ws.Range('L132').Value = 339322.665
$ws.Range('M132').Value = -15708.731
$ws.Range('N132').Value = -344382.665
$ws.Range('H136').Value = 42017.5
$ws.Range('J136').Value = 42017.5
$ws.Range('L136').Value = 126052.5
$ws.Range('N136').Value = -131152.5
$ws.Range('H137').Value = 49999.832
$ws.Range('I137').Value = 49999
$ws.Range('K137').Value = 49999
$ws.Range('M137').Value = -44899
$ws.Range('H138').Value = 49885.8
$ws.Range('J138').Value = 49885.8
$ws.Range('L138').Value = 49885.8
$ws.Range('N138').Value = -60165.8

# Sheet 7
$ws = $wb.Worksheets.Item(7)
$ws.Range('H122').Value = 4750.8184
$ws.Range('I122').Value = 4321.892
$ws.Range('J122').Value = 5632.5
$ws.Range('K122').Value = 12965.676
$ws.Range('L122').Value = 16897.5
$ws.Range('M122').Value = -10515.676
$ws.Range('N122').Value = -21797.5
$ws.Range('H133').Value = 80000
$ws.Range('J133').Value = 80000
$ws.Range('L133').Value = 80000
$ws.Range('N133').Value = -85060

# Sheet 8
$ws = $wb.Worksheets.Item(8)
$ws.Range('H126').Value = 1643.2858
$ws.Range('I126').Value = 1643.2858
$ws.Range('K126').Value = 4929.857400000001
$ws.Range('M126').Value = -2459.857400000001
$ws.Range('H132').Value = 66216.94
$ws.Range('I132').Value = 9223
$ws.Range('K132').Value = 27669
$ws.Range('M132').Value = -25139
$ws.Range('H135').Value = 105000
$ws.Range('J135').Value = 105000
$ws.Range('L135').Value = 105000
$ws.Range('N135').Value = -115140
$ws.Range('H136').Value = 46909.035
$ws.Range('I136').Value = 10947.872
$ws.Range('K136').Value = 32843.61599999999
$ws.Range('M136').Value = -30293.61599999999
$ws.Range('H137').Value = 94999
$ws.Range('J137').Value = 94999
$ws.Range('L137').Value = 94999
$ws.Range('N137').Value = -105199
$ws.Range('H138').Value = 99990
$ws.Range('J138').Value = 99990
$ws.Range('L138').Value = 99990
$ws.Range('N138').Value = -110270
$ws.Range('H141').Value = 58518
$ws.Range('J141').Value = 58518
$ws.Range('L141').Value = 58518
$ws.Range('N141').Value = -68878
